$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows to append to the bottom of the sheet (rows 30-31)
$newRows = @(
    @("183090-0", "Clio - Greek Yogurt Bar Vanilla", "1", "15.45", "15.45"),
    @("183096-7", "Clio - Greek Yogurt Bar Strawberry", "1", "15.45", "15.45")
)

$startRow = 30
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
    }
}
